$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data table (rows 2-13, columns A-T) after adding the "FAPs" sending
# cluster (Dr Hou advice) and recomputing the NATMI LR-pair statistics.
$rowsData = @(
  @(2, "ECs", "Inhbb", "Acvr2b", "ECs", 3, 1, 2.525153666666667, 7.575461000000001, 0.7378191762484796, 0.7378191762484798, 3, 1, 1.285895333333333, 3.857686, 0.3864259878905995, 0.3864259878905995, 3.247083315916222, 29.223749843246, 0.2851125040664471, 0.2851125040664471),
  @(3, "ECs", "Inhbb", "Acvr2b", "FAPs", 3, 1, 2.525153666666667, 7.575461000000001, 0.7378191762484796, 0.7378191762484798, 3, 1, 0.3517506666666667, 1.055252, 0.1057050253891921, 0.1057050253891921, 0.8882244856857779, 7.994020371172001, 0.07799119475797833, 0.07799119475797835),
  @(4, "ECs", "Inhbb", "Acvr2b", "M2", 3, 1, 2.525153666666667, 7.575461000000001, 0.7378191762484796, 0.7378191762484798, 3, 1, 0.5966156666666667, 1.789847, 0.1792897076506553, 0.1792897076506553, 1.506546238274111, 13.558916144467, 0.1322833844086372, 0.1322833844086372),
  @(5, "ECs", "Inhbb", "Acvr2b", "sCs", 3, 1, 2.525153666666667, 7.575461000000001, 0.7378191762484796, 0.7378191762484798, 3, 1, 1.093401, 3.280203, 0.328579279069553, 0.3285792790695531, 2.761005544287, 24.849049898583, 0.2424320930154169, 0.242432093015417),
  @(6, "FAPs", "Inhbb", "Acvr2b", "ECs", 2, 0.6666666666666666, 0.3739756666666667, 1.121927, 0.1092711394000877, 0.1092711394000878, 3, 1, 1.285895333333333, 3.857686, 0.3864259878905995, 0.3864259878905995, 0.480893564546889, 4.328042080922001, 0.0422252079906103, 0.04222520799061032),
  @(7, "FAPs", "Inhbb", "Acvr2b", "FAPs", 2, 0.6666666666666666, 0.3739756666666667, 1.121927, 0.1092711394000877, 0.1092711394000878, 3, 1, 0.3517506666666667, 1.055252, 0.1057050253891921, 0.1057050253891921, 0.1315461900671111, 1.183915710604, 0.01155050856459222, 0.01155050856459222),
  @(8, "FAPs", "Inhbb", "Acvr2b", "M2", 2, 0.6666666666666666, 0.3739756666666667, 1.121927, 0.1092711394000877, 0.1092711394000878, 3, 1, 0.5966156666666667, 1.789847, 0.1792897076506553, 0.1792897076506553, 0.2231197416854445, 2.008077675169, 0.01959119063769573, 0.01959119063769573),
  @(9, "FAPs", "Inhbb", "Acvr2b", "sCs", 2, 0.6666666666666666, 0.3739756666666667, 1.121927, 0.1092711394000877, 0.1092711394000878, 3, 1, 1.093401, 3.280203, 0.328579279069553, 0.3285792790695531, 0.408905367909, 3.680148311181, 0.03590423220718946, 0.03590423220718948),
  @(10, "sCs", "Inhbb", "Acvr2b", "ECs", 3, 1, 0.5233266666666666, 1.56998, 0.1529096843514326, 0.1529096843514326, 3, 1, 1.285895333333333, 3.857686, 0.3864259878905995, 0.3864259878905995, 0.6729433184755556, 6.05648986628, 0.05908827583354208, 0.05908827583354209),
  @(11, "sCs", "Inhbb", "Acvr2b", "FAPs", 3, 1, 0.5233266666666666, 1.56998, 0.1529096843514326, 0.1529096843514326, 3, 1, 0.3517506666666667, 1.055252, 0.1057050253891921, 0.1057050253891921, 0.1840805038844444, 1.65672453496, 0.01616332206662153, 0.01616332206662154),
  @(12, "sCs", "Inhbb", "Acvr2b", "M2", 3, 1, 0.5233266666666666, 1.56998, 0.1529096843514326, 0.1529096843514326, 3, 1, 0.5966156666666667, 1.789847, 0.1792897076506553, 0.1792897076506553, 0.3122248881177777, 2.81002399306, 0.02741513260432233, 0.02741513260432233),
  @(13, "sCs", "Inhbb", "Acvr2b", "sCs", 3, 1, 0.5233266666666666, 1.56998, 0.1529096843514326, 0.1529096843514326, 3, 1, 1.093401, 3.280203, 0.328579279069553, 0.3285792790695531, 0.5722059006599999, 5.149853105939999, 0.05024295384694664, 0.05024295384694665)
)

foreach ($row in $rowsData) {
  $r = $row[0]
  $ws.Cells.Item($r, 1).Value2 = $row[1]
  $ws.Cells.Item($r, 2).Value2 = $row[2]
  $ws.Cells.Item($r, 3).Value2 = $row[3]
  $ws.Cells.Item($r, 4).Value2 = $row[4]
  $ws.Cells.Item($r, 5).Value2 = $row[5]
  $ws.Cells.Item($r, 6).Value2 = $row[6]
  $ws.Cells.Item($r, 7).Value2 = $row[7]
  $ws.Cells.Item($r, 8).Value2 = $row[8]
  $ws.Cells.Item($r, 9).Value2 = $row[9]
  $ws.Cells.Item($r, 10).Value2 = $row[10]
  $ws.Cells.Item($r, 11).Value2 = $row[11]
  $ws.Cells.Item($r, 12).Value2 = $row[12]
  $ws.Cells.Item($r, 13).Value2 = $row[13]
  $ws.Cells.Item($r, 14).Value2 = $row[14]
  $ws.Cells.Item($r, 15).Value2 = $row[15]
  $ws.Cells.Item($r, 16).Value2 = $row[16]
  $ws.Cells.Item($r, 17).Value2 = $row[17]
  $ws.Cells.Item($r, 18).Value2 = $row[18]
  $ws.Cells.Item($r, 19).Value2 = $row[19]
  $ws.Cells.Item($r, 20).Value2 = $row[20]
}